$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ADS8686")
$ws.Activate()

$data = @(
    ,("config", "0x02", 0, 8, 7, 0)
    ,("chan_sel", "0x03", 0, 8, 7, 0)
    ,("rangeA1", "0x04", "0xFF", 8, 7, 0)
    ,("rangeA2", "0x05", "0xFF", 8, 7, 0)
    ,("rangeB1", "0x06", "0xFF", 8, 7, 0)
    ,("rangeB2", "0x07", "0xFF", 8, 7, 0)
    ,("status", "0x08", "0xFF", 16, 15, 0)
    ,("over_rangeA", "0xA", "0x00", 8, 7, 0)
    ,("over_rangeB", "0xB", "0x00", 8, 7, 0)
    ,("lpf", "0xC", 0, 2, 1, 0)
    ,("devID", "0x10", 2, 2, 1, 0)
    ,("seq0", "0x20", 0, 8, 7, 0)
    ,("seq1", "0x21", "0x11", 8, 7, 0)
    ,("seq2", "0x22", "0x22", 8, 7, 0)
    ,("seq3", "0x23", "0x33", 8, 7, 0)
    ,("seq4", "0x24", "0x44", 8, 7, 0)
    ,("seq5", "0x25", "0x55", 8, 7, 0)
    ,("seq6", "0x26", "0x66", 8, 7, 0)
    ,("seq7", "0x27", "0x77", 8, 7, 0)
    ,("seq8", "0x28", 0, 8, 7, 0)
    ,("seq9", "0x29", 0, 8, 7, 0)
    ,("seq10", "0x2A", 0, 8, 7, 0)
    ,("seq11", "0x2B", 0, 8, 7, 0)
    ,("seq12", "0x2C", 0, 8, 7, 0)
    ,("seq13", "0x2D", 0, 8, 7, 0)
    ,("seq14", "0x2E", 0, 8, 7, 0)
    ,("seq15", "0x2F", 0, 8, 7, 0)
    ,("seq16", "0x30", 0, 8, 7, 0)
    ,("seq17", "0x31", 0, 8, 7, 0)
    ,("seq18", "0x32", 0, 8, 7, 0)
    ,("seq19", "0x33", 0, 8, 7, 0)
    ,("seq20", "0x34", 0, 8, 7, 0)
    ,("seq21", "0x35", 0, 8, 7, 0)
    ,("seq22", "0x36", 0, 8, 7, 0)
    ,("seq23", "0x37", 0, 8, 7, 0)
    ,("seq24", "0x38", 0, 8, 7, 0)
    ,("seq25", "0x39", 0, 8, 7, 0)
    ,("seq26", "0x3A", 0, 8, 7, 0)
    ,("seq27", "0x3B", 0, 8, 7, 0)
    ,("seq28", "0x3C", 0, 8, 7, 0)
    ,("seq29", "0x3D", 0, 8, 7, 0)
    ,("seq30", "0x3E", 0, 8, 7, 0)
    ,("seq31", "0x3F", 0, 8, 7, 0)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

$ws.Range("A2:F44").Select()
